{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2024-10-20 Sunday\", \"2024-10-21 Monday\"],\n  [\"243\u00f74=\", \"779\u00f79=\"],\n  [\"361\u00f75=\", \"724\u00f78=\"],\n  [\"655\u00f73=\", \"245\u00f76=\"],\n  [\"971\u00f75=\", \"630\u00f77=\"],\n  [\"895\u00f79=\", \"553\u00f79=\"],\n  [\"905\u00f72=\", \"859\u00f72=\"],\n  [\"170\u00f72=\", \"135\u00f76=\"],\n  [\"665\u00f78=\", \"916\u00f79=\"],\n  [\"638\u00f77=\", \"403\u00f77=\"],\n  [\"764\u00f73=\", \"941\u00f78=\"],\n  [\"489\u00f76=\", \"756\u00f73=\"],\n  [\"772\u00f79=\", \"316\u00f77=\"],\n  [\"834\u00f74=\", \"224\u00f79=\"],\n  [\"386\u00f77=\", \"347\u00f76=\"],\n  [\"526\u00f79=\", \"827\u00f75=\"],\n  [\"981\u00f73=\", \"146\u00f74=\"],\n  [\"324\u00f72=\", \"640\u00f78=\"],\n  [\"238\u00f72=\", \"734\u00f73=\"],\n  [\"686\u00f74=\", \"577\u00f73=\"],\n  [\"535\u00f73=\", \"700\u00f78=\"],\n  [\"153\u00f79=\", \"640\u00f78=\"],\n  [\"661\u00f76=\", \"100\u00f73=\"],\n  [\"878\u00f72=\", \"901\u00f74=\"],\n  [\"828\u00f78=\", \"689\u00f72=\"],\n  [\"996\u00f72=\", \"106\u00f77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-10-20 Sunday\", \"2024-10-21 Monday\")\n    ,@(\"243\u00f74=\", \"779\u00f79=\")\n    ,@(\"361\u00f75=\", \"724\u00f78=\")\n    ,@(\"655\u00f73=\", \"245\u00f76=\")\n    ,@(\"971\u00f75=\", \"630\u00f77=\")\n    ,@(\"895\u00f79=\", \"553\u00f79=\")\n    ,@(\"905\u00f72=\", \"859\u00f72=\")\n    ,@(\"170\u00f72=\", \"135\u00f76=\")\n    ,@(\"665\u00f78=\", \"916\u00f79=\")\n    ,@(\"638\u00f77=\", \"403\u00f77=\")\n    ,@(\"764\u00f73=\", \"941\u00f78=\")\n    ,@(\"489\u00f76=\", \"756\u00f73=\")\n    ,@(\"772\u00f79=\", \"316\u00f77=\")\n    ,@(\"834\u00f74=\", \"224\u00f79=\")\n    ,@(\"386\u00f77=\", \"347\u00f76=\")\n    ,@(\"526\u00f79=\", \"827\u00f75=\")\n    ,@(\"981\u00f73=\", \"146\u00f74=\")\n    ,@(\"324\u00f72=\", \"640\u00f78=\")\n    ,@(\"238\u00f72=\", \"734\u00f73=\")\n    ,@(\"686\u00f74=\", \"577\u00f73=\")\n    ,@(\"535\u00f73=\", \"700\u00f78=\")\n    ,@(\"153\u00f79=\", \"640\u00f78=\")\n    ,@(\"661\u00f76=\", \"100\u00f73=\")\n    ,@(\"878\u00f72=\", \"901\u00f74=\")\n    ,@(\"828\u00f78=\", \"689\u00f72=\")\n    ,@(\"996\u00f72=\", \"106\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
